# EDAC-to-OBOE.xlsx -- "updated mappings and pres"
#
# Updates a few mapping values on Sheet1 and leaves the selection where the
# author last left it (E12), plus a best-effort nudge of the saved window
# size to match the slightly taller Excel window recorded in the workbook
# view the last time the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated mapping values (column D, "FGDC attrunit"/count column) -----
$ws.Range("D2").Value = 1532
$ws.Range("D3").Value = 1532
$ws.Range("D4").Value = 3064

# --- Window chrome: the author's Excel window grew a touch taller ---------
$win = $excel.ActiveWindow
$win.Height = 16480

# --- Leave the selection where the author left it before saving -----------
$ws.Range("E12").Select()
